$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Target cluster = ECs (same cluster), updated stats
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf9"
$ws.Range("C2").Value = "Fgfr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.883656666666667
$ws.Range("H2").Value = 5.650970000000001
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3341516666666666
$ws.Range("N2").Value = 1.002455
$ws.Range("O2").Value = 0.02354483703663662
$ws.Range("P2").Value = 0.03467746452042096
$ws.Range("Q2").Value = 0.6294270145944445
$ws.Range("R2").Value = 5.664843131350001
$ws.Range("S2").Value = 0.02354483703663662
$ws.Range("T2").Value = 0.03467746452042096

# Row 3: Target cluster = FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf9"
$ws.Range("C3").Value = "Fgfr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.883656666666667
$ws.Range("H3").Value = 5.650970000000001
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.153566
$ws.Range("N3").Value = 0.460698
$ws.Range("O3").Value = 0.01082049501783563
$ws.Range("P3").Value = 0.01593671391696275
$ws.Range("Q3").Value = 0.2892656196733334
$ws.Range("R3").Value = 2.60339057706
$ws.Range("S3").Value = 0.01082049501783563
$ws.Range("T3").Value = 0.01593671391696275

# Row 4: Target cluster = M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf9"
$ws.Range("C4").Value = "Fgfr4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.883656666666667
$ws.Range("H4").Value = 5.650970000000001
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03596566666666667
$ws.Range("N4").Value = 0.107897
$ws.Range("O4").Value = 0.002534195830976933
$ws.Range("P4").Value = 0.003732431270590561
$ws.Range("Q4").Value = 0.06774696778777779
$ws.Range("R4").Value = 0.6097227100900001
$ws.Range("S4").Value = 0.002534195830976933
$ws.Range("T4").Value = 0.003732431270590561

# Row 5: Target cluster = sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf9"
$ws.Range("C5").Value = "Fgfr4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.883656666666667
$ws.Range("H5").Value = 5.650970000000001
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.6684585
$ws.Range("N5").Value = 27.336917
$ws.Range("O5").Value = 0.9631004721145509
$ws.Range("P5").Value = 0.9456533902920258
$ws.Range("Q5").Value = 25.74668297658167
$ws.Range("R5").Value = 154.48009785949
$ws.Range("S5").Value = 0.9631004721145509
$ws.Range("T5").Value = 0.9456533902920258
